# Cascade workbook: add a new "Cascade Characteristics" worksheet between
# "Compartments" and "Transitions", matching the author's commit.

$wb = $excel.ActiveWorkbook

$compartments = $wb.Worksheets.Item("Compartments")
$transitions  = $wb.Worksheets.Item("Transitions")

# Insert the new worksheet right after "Compartments" (i.e. before "Transitions")
$new = $wb.Worksheets.Add($null, $compartments)
$new.Name = "Cascade Characteristics"

# --- Header row ---
$new.Cells.Item(1,1).Value = 'Code Label'
$new.Cells.Item(1,2).Value = 'Full Name'
$new.Cells.Item(1,3).Value = 'Denominator'
$new.Cells.Item(1,4).Value = 'Includes'
$new.Range("A1:D1").Font.Bold = $true

# --- Data rows ---
$new.Cells.Item(2,1).Value = 'lt_inf'
$new.Cells.Item(2,2).Value = 'Latent Infections'
$new.Cells.Item(2,4).Value = 'lte'
$new.Cells.Item(2,5).Value = 'ltsu'
$new.Cells.Item(2,6).Value = 'ltsd'
$new.Cells.Item(2,7).Value = 'ltst'
$new.Cells.Item(2,8).Value = 'ltfu'
$new.Cells.Item(2,9).Value = 'ltfd'
$new.Cells.Item(2,10).Value = 'ltft'

$new.Cells.Item(3,1).Value = 's+_inf'
$new.Cells.Item(3,2).Value = 'Smear Positive Infections'
$new.Cells.Item(3,4).Value = 's+e'
$new.Cells.Item(3,5).Value = 's+du'
$new.Cells.Item(3,6).Value = 's+dd'
$new.Cells.Item(3,7).Value = 's+dt'
$new.Cells.Item(3,8).Value = 's+mu'
$new.Cells.Item(3,9).Value = 's+md'
$new.Cells.Item(3,10).Value = 's+mt'
$new.Cells.Item(3,11).Value = 's+xu'
$new.Cells.Item(3,12).Value = 's+xd'
$new.Cells.Item(3,13).Value = 's+xt'

$new.Cells.Item(4,1).Value = 's-_inf'
$new.Cells.Item(4,2).Value = 'Smear Negative Infections'
$new.Cells.Item(4,4).Value = 's-e'
$new.Cells.Item(4,5).Value = 's-du'
$new.Cells.Item(4,6).Value = 's-dd'
$new.Cells.Item(4,7).Value = 's-dt'
$new.Cells.Item(4,8).Value = 's-mu'
$new.Cells.Item(4,9).Value = 's-md'
$new.Cells.Item(4,10).Value = 's-mt'
$new.Cells.Item(4,11).Value = 's-xu'
$new.Cells.Item(4,12).Value = 's-xd'
$new.Cells.Item(4,13).Value = 's-xt'

$new.Cells.Item(5,1).Value = 'ac_inf'
$new.Cells.Item(5,2).Value = 'Active Infections'
$new.Cells.Item(5,4).Value = 's+inf'
$new.Cells.Item(5,5).Value = 's-_inf'

$new.Cells.Item(6,1).Value = 'alive'
$new.Cells.Item(6,2).Value = 'Total Living Population'
$new.Cells.Item(6,4).Value = 'sus'
$new.Cells.Item(6,5).Value = 'vac'
$new.Cells.Item(6,6).Value = 'rec'
$new.Cells.Item(6,7).Value = 'lt_inf'
$new.Cells.Item(6,8).Value = 'ac_inf'

$new.Cells.Item(7,1).Value = 'lt_prev'
$new.Cells.Item(7,2).Value = 'Latent Prevalence'
$new.Cells.Item(7,3).Value = 'alive'
$new.Cells.Item(7,4).Value = 'lt_inf'

$new.Cells.Item(8,1).Value = 's+_prev'
$new.Cells.Item(8,2).Value = 'Smear Positive Prevalence'
$new.Cells.Item(8,3).Value = 'alive'
$new.Cells.Item(8,4).Value = 's+_inf'

$new.Cells.Item(9,1).Value = 's-_prev'
$new.Cells.Item(9,2).Value = 'Smear Negative Prevalence'
$new.Cells.Item(9,3).Value = 'alive'
$new.Cells.Item(9,4).Value = 's-_inf'

$new.Cells.Item(10,1).Value = 'ac_prev'
$new.Cells.Item(10,2).Value = 'Active Prevalence'
$new.Cells.Item(10,3).Value = 'alive'
$new.Cells.Item(10,4).Value = 'ac_inf'

# Approximate column widths (best-effort, cosmetic only)
$widths = @(10.21875,23.44140625,12.109375,8,5.5546875,4.77734375,5.21875,6,5.33203125,5,4.6640625,4.6640625,4.33203125)
for ($i = 1; $i -le $widths.Length; $i++) {
    $new.Columns.Item($i).ColumnWidth = $widths[$i-1]
}

$new.PageSetup.Orientation = 1

# --- View / selection state ---

# "Compartments" keeps selection A2:A31 and scrolls down a couple of rows,
# and is no longer the active/selected tab.
$compartments.Range("A2:A31").Select() | Out-Null

# The new sheet becomes the active tab with selection on E15.
$new.Activate() | Out-Null
$new.Range("E15").Select() | Out-Null

Write-Host "Cascade Characteristics worksheet added."
